# Edit script for "Topic 10 Generalized linear mixed effects models.pptx"
# Applies:
#   1) Date placeholder text update (2018-02-23 -> 2018-03-01) on the slide
#      master and every slide layout.
#   2) Repositioning / re-alignment of the Title, Picture and TextBox shapes
#      on slide 1 (the title slide).
#
# NOTE: Shape.Left/.Top/.Width/.Height are expressed in points (1 pt =
# 12700 EMU) and are stored as single-precision floats by the host, which
# truncates when converting back to EMU. Adding half an EMU before the
# division keeps the round-trip exact for the target EMU values below.

function EMU([double]$emu) {
    return ($emu + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation

# --- 1) Update the cached "datetimeFigureOut" field text everywhere it
#        appears: the slide master and all custom (slide) layouts. ---
$master = $p.SlideMaster

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "2018-03-01"
        }
    }
}

Update-DatePlaceholder $master.Shapes

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}

# --- 2) Slide 1 shape tweaks ---
$s1 = $p.Slides.Item(1)

# Title 1 : move up slightly and center the title text
$title = $s1.Shapes.Item(1)
$title.Top = EMU(980728)
$title.TextFrame.TextRange.ParagraphFormat.Alignment = 2   # ppAlignCenter

# Picture 2 : reposition / resize the picture
$pic = $s1.Shapes.Item(2)
$pic.Left = EMU(3491880)
$pic.Top = EMU(2924944)
$pic.Width = EMU(2179344)
$pic.Height = EMU(3280022)

# TextBox 3 : reposition the caption textbox
$tbox = $s1.Shapes.Item(3)
$tbox.Left = EMU(3988551)
$tbox.Top = EMU(6309320)
